# Fix locus naming convention: "Ots37124-XXXXXXX" -> "Ots37124.XXXXXXX"
# (hyphen changed to dot) on both the "ALL" and "diag" worksheets, and
# restore "ALL" as the active/selected sheet & cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ALL")
$ws2 = $wb.Worksheets.Item("diag")

# --- "ALL" sheet: rows 80-94, column A (3 rows per locus) ---
$ws1.Range("A80").Value = "Ots37124.12267397"
$ws1.Range("A81").Value = "Ots37124.12267397"
$ws1.Range("A82").Value = "Ots37124.12267397"

$ws1.Range("A83").Value = "Ots37124.12270118"
$ws1.Range("A84").Value = "Ots37124.12270118"
$ws1.Range("A85").Value = "Ots37124.12270118"

$ws1.Range("A86").Value = "Ots37124.12272852"
$ws1.Range("A87").Value = "Ots37124.12272852"
$ws1.Range("A88").Value = "Ots37124.12272852"

$ws1.Range("A89").Value = "Ots37124.12277401"
$ws1.Range("A90").Value = "Ots37124.12277401"
$ws1.Range("A91").Value = "Ots37124.12277401"

$ws1.Range("A92").Value = "Ots37124.12310649"
$ws1.Range("A93").Value = "Ots37124.12310649"
$ws1.Range("A94").Value = "Ots37124.12310649"

# --- "diag" sheet: rows 28-32, column A (one row per locus) ---
$ws2.Range("A28").Value = "Ots37124.12267397"
$ws2.Range("A29").Value = "Ots37124.12270118"
$ws2.Range("A30").Value = "Ots37124.12272852"
$ws2.Range("A31").Value = "Ots37124.12277401"
$ws2.Range("A32").Value = "Ots37124.12310649"

# --- restore view state: "diag" selection moves to A32 (no longer the ---
# --- active tab), and "ALL" becomes the active tab with Q12 selected. ---
$null = $ws2.Range("A32").Select()
$null = $ws1.Activate()
$null = $ws1.Range("Q12").Select()
